# Add "Unsuitable exception type" check to the Workflow checklist sheet.
# This appends a new row (32) below the last existing check row (31),
# reusing that row's formatting and filling in the new check's details.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$srcRow = 31
$newRow = 32

# Copy formatting (styles) from the row above so the new row matches
# the look of the rest of the checklist table.
$ws.Range("A$srcRow`:G$srcRow").Copy()
$ws.Range("A$newRow`:G$newRow").PasteSpecial(-4122)
$ws.Rows.Item($newRow).RowHeight = 85

# Fill in the new check's data.
$ws.Range("A$newRow").Value = "No"
$ws.Range("B$newRow").Value = "Unsuitable exception type"
$ws.Range("C$newRow").Value = "Checks\Custom\UnsuitableExceptionType.xaml"
$ws.Range("E$newRow").Value = "Fix"
$ws.Range("F$newRow").Value = "When throwing exceptions, it is recommended to properly distinguish between application-originated and business-originated exceptions. The type of the exception to be thrown or caught should be as specific as possible, and Exception and ApplicationException should be avoided."
$ws.Range("G$newRow").Value = "Use specific exception types and avoid using generic types such as Exception and ApplicationException."
